# Update the "K" column (column G) values for rows 2-33 on the active sheet.
# These values were regenerated upstream (K replaces the old Strike# figures);
# here we just write the final computed values into the corresponding cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 4
    4  = 3
    5  = 6
    6  = 10
    7  = 7
    8  = 5
    9  = 7
    10 = 7
    11 = 5
    12 = 5
    13 = 6
    14 = 5
    15 = 6
    16 = 1
    17 = 4
    18 = 7
    19 = 3
    20 = 3
    21 = 5
    22 = 10
    23 = 9
    24 = 5
    25 = 6
    26 = 9
    27 = 4
    28 = 5
    29 = 2
    30 = 8
    31 = 5
    32 = 2
    33 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
